$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.123.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.912.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.688.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("E14").Value = "  +2.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.091.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0739"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.536.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.607"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.938"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.35%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("E39").Value = "  -1.16%  "

$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.85%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.821.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "
